$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.370286
$ws.Range("H2").Value = 7.110858
$ws.Range("I2").Value = 0.3026841782318013
$ws.Range("J2").Value = 0.3026841782318014
$ws.Range("O2").Value = 0.7091726973716084
$ws.Range("P2").Value = 0.7091726973716084
$ws.Range("Q2").Value = 4.047942826986667
$ws.Range("R2").Value = 36.43148544288
$ws.Range("S2").Value = 0.2146553551283552
$ws.Range("T2").Value = 0.2146553551283553

# Row 3
$ws.Range("G3").Value = 2.370286
$ws.Range("H3").Value = 7.110858
$ws.Range("I3").Value = 0.3026841782318013
$ws.Range("J3").Value = 0.3026841782318014
$ws.Range("M3").Value = 0.7003526666666667
$ws.Range("N3").Value = 2.101058
$ws.Range("O3").Value = 0.2908273026283917
$ws.Range("P3").Value = 0.2908273026283917
$ws.Range("Q3").Value = 1.660036120862667
$ws.Range("R3").Value = 14.940325087764
$ws.Range("S3").Value = 0.08802882310344613
$ws.Range("T3").Value = 0.08802882310344615

# Row 4
$ws.Range("I4").Value = 0.2022126055089961
$ws.Range("J4").Value = 0.2022126055089961
$ws.Range("O4").Value = 0.7091726973716084
$ws.Range("P4").Value = 0.7091726973716084
$ws.Range("S4").Value = 0.1434036588913557
$ws.Range("T4").Value = 0.1434036588913558

# Row 5
$ws.Range("I5").Value = 0.2022126055089961
$ws.Range("J5").Value = 0.2022126055089961
$ws.Range("M5").Value = 0.7003526666666667
$ws.Range("N5").Value = 2.101058
$ws.Range("O5").Value = 0.2908273026283917
$ws.Range("P5").Value = 0.2908273026283917
$ws.Range("Q5").Value = 1.109011482528222
$ws.Range("R5").Value = 9.981103342754
$ws.Range("S5").Value = 0.0588089466176404
$ws.Range("T5").Value = 0.0588089466176404

# Row 6
$ws.Range("G6").Value = 2.286703333333333
$ws.Range("H6").Value = 6.860109999999999
$ws.Range("I6").Value = 0.2920107190904054
$ws.Range("J6").Value = 0.2920107190904054
$ws.Range("O6").Value = 0.7091726973716084
$ws.Range("P6").Value = 0.7091726973716084
$ws.Range("Q6").Value = 3.905201463288888
$ws.Range("R6").Value = 35.14681316959999
$ws.Range("S6").Value = 0.2070860293187658
$ws.Range("T6").Value = 0.2070860293187658

# Row 7
$ws.Range("G7").Value = 2.286703333333333
$ws.Range("H7").Value = 6.860109999999999
$ws.Range("I7").Value = 0.2920107190904054
$ws.Range("J7").Value = 0.2920107190904054
$ws.Range("M7").Value = 0.7003526666666667
$ws.Range("N7").Value = 2.101058
$ws.Range("O7").Value = 0.2908273026283917
$ws.Range("P7").Value = 0.2908273026283917
$ws.Range("Q7").Value = 1.601498777375556
$ws.Range("R7").Value = 14.41348899638
$ws.Range("S7").Value = 0.08492468977163961
$ws.Range("T7").Value = 0.08492468977163961

# Row 8
$ws.Range("G8").Value = 1.590394666666667
$ws.Range("H8").Value = 4.771184
$ws.Range("I8").Value = 0.2030924971687972
$ws.Range("J8").Value = 0.2030924971687972
$ws.Range("O8").Value = 0.7091726973716084
$ws.Range("P8").Value = 0.7091726973716084
$ws.Range("Q8").Value = 2.716054806471111
$ws.Range("R8").Value = 24.44449325824
$ws.Range("S8").Value = 0.1440276540331316
$ws.Range("T8").Value = 0.1440276540331316

# Row 9
$ws.Range("G9").Value = 1.590394666666667
$ws.Range("H9").Value = 4.771184
$ws.Range("I9").Value = 0.2030924971687972
$ws.Range("J9").Value = 0.2030924971687972
$ws.Range("M9").Value = 0.7003526666666667
$ws.Range("N9").Value = 2.101058
$ws.Range("O9").Value = 0.2908273026283917
$ws.Range("P9").Value = 0.2908273026283917
$ws.Range("Q9").Value = 1.113837145852445
$ws.Range("R9").Value = 10.024534312672
$ws.Range("S9").Value = 0.05906484313566556
$ws.Range("T9").Value = 0.05906484313566556
